$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Range("A32").Value = 112211016
$ws.Range("B32").Value = 90666
$ws.Range("C32").Value = "Ovaliderad"
$ws.Range("D32").Value = "LC"
$ws.Range("E32").Value = 4364
$ws.Range("F32").Value = "Dropptaggsvamp"
$ws.Range("G32").Value = "Hydnellum ferrugineum"
$ws.Range("H32").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I32").NumberFormat = "@"
$ws.Range("I32").Value = "1"
$ws.Range("J32").Value = "fruktkroppar"
$ws.Range("K32").NumberFormat = "@"
$ws.Range("K32").Value = ""
$ws.Range("P32").Value = "S Simsbodarna, Dlr"
$ws.Range("Q32").Value = 515070
$ws.Range("R32").Value = 6704834
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = "Dalarna"
$ws.Range("U32").Value = "Borlänge"
$ws.Range("V32").Value = "Dalarna"
$ws.Range("W32").Value = "Stora Tuna"
$ws.Range("Y32").NumberFormat = "@"
$ws.Range("Y32").Value = "2023-09-20"
$ws.Range("Z32").Value = "10:55"
$ws.Range("AA32").NumberFormat = "@"
$ws.Range("AA32").Value = "2023-09-20"
$ws.Range("AB32").Value = "10:55"
$ws.Range("AD32").Value = $false
$ws.Range("AE32").Value = $false
$ws.Range("AG32").Value = $false
$ws.Range("AT32").NumberFormat = "@"
$ws.Range("AT32").Value = ""
$ws.Range("AW32").Value = "Lars-Erik Nilsson"
$ws.Range("AX32").Value = "Lars-Erik Nilsson"
$ws.Range("AY32").NumberFormat = "@"
$ws.Range("AY32").Value = ""

# Row 33
$ws.Range("A33").Value = 112211348
$ws.Range("B33").Value = 89369
$ws.Range("C33").Value = "Ovaliderad"
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 5447
$ws.Range("F33").Value = "Vedticka"
$ws.Range("G33").Value = "Fuscoporia viticola"
$ws.Range("H33").Value = "(Schwein.) Murrill"
$ws.Range("I33").NumberFormat = "@"
$ws.Range("I33").Value = ""
$ws.Range("K33").NumberFormat = "@"
$ws.Range("K33").Value = ""
$ws.Range("P33").Value = "Simsbodarna S, Dlr"
$ws.Range("Q33").Value = 515173
$ws.Range("R33").Value = 6704768
$ws.Range("S33").Value = 1
$ws.Range("T33").Value = "Dalarna"
$ws.Range("U33").Value = "Borlänge"
$ws.Range("V33").Value = "Dalarna"
$ws.Range("W33").Value = "Stora Tuna"
$ws.Range("Y33").NumberFormat = "@"
$ws.Range("Y33").Value = "2023-09-20"
$ws.Range("Z33").Value = "11:15"
$ws.Range("AA33").NumberFormat = "@"
$ws.Range("AA33").Value = "2023-09-20"
$ws.Range("AB33").Value = "11:15"
$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AG33").Value = $false
$ws.Range("AT33").NumberFormat = "@"
$ws.Range("AT33").Value = ""
$ws.Range("AW33").Value = "Lars-Erik Nilsson"
$ws.Range("AX33").Value = "Lars-Erik Nilsson"
$ws.Range("AY33").NumberFormat = "@"
$ws.Range("AY33").Value = ""

# Row 34
$ws.Range("A34").Value = 112212788
$ws.Range("B34").Value = 90666
$ws.Range("C34").Value = "Ovaliderad"
$ws.Range("D34").Value = "LC"
$ws.Range("E34").Value = 4364
$ws.Range("F34").Value = "Dropptaggsvamp"
$ws.Range("G34").Value = "Hydnellum ferrugineum"
$ws.Range("H34").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I34").NumberFormat = "@"
$ws.Range("I34").Value = "2"
$ws.Range("J34").Value = "fruktkroppar"
$ws.Range("K34").NumberFormat = "@"
$ws.Range("K34").Value = ""
$ws.Range("P34").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q34").Value = 515101
$ws.Range("R34").Value = 6704641
$ws.Range("S34").Value = 1
$ws.Range("T34").Value = "Dalarna"
$ws.Range("U34").Value = "Borlänge"
$ws.Range("V34").Value = "Dalarna"
$ws.Range("W34").Value = "Stora Tuna"
$ws.Range("Y34").NumberFormat = "@"
$ws.Range("Y34").Value = "2023-09-20"
$ws.Range("Z34").Value = "12:37"
$ws.Range("AA34").NumberFormat = "@"
$ws.Range("AA34").Value = "2023-09-20"
$ws.Range("AB34").Value = "12:37"
$ws.Range("AD34").Value = $false
$ws.Range("AE34").Value = $false
$ws.Range("AG34").Value = $false
$ws.Range("AT34").NumberFormat = "@"
$ws.Range("AT34").Value = ""
$ws.Range("AW34").Value = "Lars-Erik Nilsson"
$ws.Range("AX34").Value = "Lars-Erik Nilsson"
$ws.Range("AY34").NumberFormat = "@"
$ws.Range("AY34").Value = ""

# Row 35
$ws.Range("A35").Value = 112212094
$ws.Range("B35").Value = 90666
$ws.Range("C35").Value = "Ovaliderad"
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 4364
$ws.Range("F35").Value = "Dropptaggsvamp"
$ws.Range("G35").Value = "Hydnellum ferrugineum"
$ws.Range("H35").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I35").NumberFormat = "@"
$ws.Range("I35").Value = "10"
$ws.Range("J35").Value = "fruktkroppar"
$ws.Range("K35").NumberFormat = "@"
$ws.Range("K35").Value = ""
$ws.Range("P35").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q35").Value = 515450
$ws.Range("R35").Value = 6704585
$ws.Range("S35").Value = 2
$ws.Range("T35").Value = "Dalarna"
$ws.Range("U35").Value = "Borlänge"
$ws.Range("V35").Value = "Dalarna"
$ws.Range("W35").Value = "Stora Tuna"
$ws.Range("Y35").NumberFormat = "@"
$ws.Range("Y35").Value = "2023-09-20"
$ws.Range("Z35").Value = "11:59"
$ws.Range("AA35").NumberFormat = "@"
$ws.Range("AA35").Value = "2023-09-20"
$ws.Range("AB35").Value = "11:59"
$ws.Range("AC35").Value = "Ca. Färska fk."
$ws.Range("AD35").Value = $false
$ws.Range("AE35").Value = $false
$ws.Range("AG35").Value = $false
$ws.Range("AT35").NumberFormat = "@"
$ws.Range("AT35").Value = ""
$ws.Range("AW35").Value = "Lars-Erik Nilsson"
$ws.Range("AX35").Value = "Lars-Erik Nilsson"
$ws.Range("AY35").NumberFormat = "@"
$ws.Range("AY35").Value = ""

# Row 36
$ws.Range("A36").Value = 112212237
$ws.Range("B36").Value = 90689
$ws.Range("C36").Value = "Ovaliderad"
$ws.Range("D36").Value = "NT"
$ws.Range("E36").Value = 5966
$ws.Range("F36").Value = "Motaggsvamp"
$ws.Range("G36").Value = "Sarcodon squamosus"
$ws.Range("H36").Value = "(Schaeff.) Quél."
$ws.Range("I36").NumberFormat = "@"
$ws.Range("I36").Value = "2"
$ws.Range("J36").Value = "fruktkroppar"
$ws.Range("K36").NumberFormat = "@"
$ws.Range("K36").Value = ""
$ws.Range("P36").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q36").Value = 515492
$ws.Range("R36").Value = 6704591
$ws.Range("S36").Value = 1
$ws.Range("T36").Value = "Dalarna"
$ws.Range("U36").Value = "Borlänge"
$ws.Range("V36").Value = "Dalarna"
$ws.Range("W36").Value = "Stora Tuna"
$ws.Range("Y36").NumberFormat = "@"
$ws.Range("Y36").Value = "2023-09-20"
$ws.Range("Z36").Value = "11:59"
$ws.Range("AA36").NumberFormat = "@"
$ws.Range("AA36").Value = "2023-09-20"
$ws.Range("AB36").Value = "11:59"
$ws.Range("AC36").Value = "Stora fina ex."
$ws.Range("AD36").Value = $false
$ws.Range("AE36").Value = $false
$ws.Range("AG36").Value = $false
$ws.Range("AT36").NumberFormat = "@"
$ws.Range("AT36").Value = ""
$ws.Range("AW36").Value = "Lars-Erik Nilsson"
$ws.Range("AX36").Value = "Lars-Erik Nilsson"
$ws.Range("AY36").NumberFormat = "@"
$ws.Range("AY36").Value = ""

# Row 37
$ws.Range("A37").Value = 112211929
$ws.Range("B37").Value = 90666
$ws.Range("C37").Value = "Ovaliderad"
$ws.Range("D37").Value = "LC"
$ws.Range("E37").Value = 4364
$ws.Range("F37").Value = "Dropptaggsvamp"
$ws.Range("G37").Value = "Hydnellum ferrugineum"
$ws.Range("H37").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I37").NumberFormat = "@"
$ws.Range("I37").Value = "10"
$ws.Range("J37").Value = "fruktkroppar"
$ws.Range("K37").NumberFormat = "@"
$ws.Range("K37").Value = ""
$ws.Range("P37").Value = "SO Simsbodarna, Dlr"
$ws.Range("Q37").Value = 515370
$ws.Range("R37").Value = 6704610
$ws.Range("S37").Value = 5
$ws.Range("T37").Value = "Dalarna"
$ws.Range("U37").Value = "Borlänge"
$ws.Range("V37").Value = "Dalarna"
$ws.Range("W37").Value = "Stora Tuna"
$ws.Range("Y37").NumberFormat = "@"
$ws.Range("Y37").Value = "2023-09-20"
$ws.Range("Z37").Value = "11:45"
$ws.Range("AA37").NumberFormat = "@"
$ws.Range("AA37").Value = "2023-09-20"
$ws.Range("AB37").Value = "11:45"
$ws.Range("AC37").Value = "G:a fk."
$ws.Range("AD37").Value = $false
$ws.Range("AE37").Value = $false
$ws.Range("AG37").Value = $false
$ws.Range("AT37").NumberFormat = "@"
$ws.Range("AT37").Value = ""
$ws.Range("AW37").Value = "Lars-Erik Nilsson"
$ws.Range("AX37").Value = "Lars-Erik Nilsson"
$ws.Range("AY37").NumberFormat = "@"
$ws.Range("AY37").Value = ""

# Row 38
$ws.Range("A38").Value = 112212369
$ws.Range("B38").Value = 90666
$ws.Range("C38").Value = "Ovaliderad"
$ws.Range("D38").Value = "LC"
$ws.Range("E38").Value = 4364
$ws.Range("F38").Value = "Dropptaggsvamp"
$ws.Range("G38").Value = "Hydnellum ferrugineum"
$ws.Range("H38").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I38").NumberFormat = "@"
$ws.Range("I38").Value = "10"
$ws.Range("J38").Value = "fruktkroppar"
$ws.Range("K38").NumberFormat = "@"
$ws.Range("K38").Value = ""
$ws.Range("P38").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q38").Value = 515497
$ws.Range("R38").Value = 6704555
$ws.Range("S38").Value = 2
$ws.Range("T38").Value = "Dalarna"
$ws.Range("U38").Value = "Borlänge"
$ws.Range("V38").Value = "Dalarna"
$ws.Range("W38").Value = "Stora Tuna"
$ws.Range("Y38").NumberFormat = "@"
$ws.Range("Y38").Value = "2023-09-20"
$ws.Range("Z38").Value = "12:12"
$ws.Range("AA38").NumberFormat = "@"
$ws.Range("AA38").Value = "2023-09-20"
$ws.Range("AB38").Value = "12:12"
$ws.Range("AD38").Value = $false
$ws.Range("AE38").Value = $false
$ws.Range("AG38").Value = $false
$ws.Range("AT38").NumberFormat = "@"
$ws.Range("AT38").Value = ""
$ws.Range("AW38").Value = "Lars-Erik Nilsson"
$ws.Range("AX38").Value = "Lars-Erik Nilsson"
$ws.Range("AY38").NumberFormat = "@"
$ws.Range("AY38").Value = ""

# Row 39
$ws.Range("A39").Value = 112212309
$ws.Range("B39").Value = 90666
$ws.Range("C39").Value = "Ovaliderad"
$ws.Range("D39").Value = "LC"
$ws.Range("E39").Value = 4364
$ws.Range("F39").Value = "Dropptaggsvamp"
$ws.Range("G39").Value = "Hydnellum ferrugineum"
$ws.Range("H39").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I39").NumberFormat = "@"
$ws.Range("I39").Value = "2"
$ws.Range("J39").Value = "fruktkroppar"
$ws.Range("K39").NumberFormat = "@"
$ws.Range("K39").Value = ""
$ws.Range("P39").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q39").Value = 515492
$ws.Range("R39").Value = 6704591
$ws.Range("S39").Value = 1
$ws.Range("T39").Value = "Dalarna"
$ws.Range("U39").Value = "Borlänge"
$ws.Range("V39").Value = "Dalarna"
$ws.Range("W39").Value = "Stora Tuna"
$ws.Range("Y39").NumberFormat = "@"
$ws.Range("Y39").Value = "2023-09-20"
$ws.Range("Z39").Value = "12:08"
$ws.Range("AA39").NumberFormat = "@"
$ws.Range("AA39").Value = "2023-09-20"
$ws.Range("AB39").Value = "12:08"
$ws.Range("AD39").Value = $false
$ws.Range("AE39").Value = $false
$ws.Range("AG39").Value = $false
$ws.Range("AT39").NumberFormat = "@"
$ws.Range("AT39").Value = ""
$ws.Range("AW39").Value = "Lars-Erik Nilsson"
$ws.Range("AX39").Value = "Lars-Erik Nilsson"
$ws.Range("AY39").NumberFormat = "@"
$ws.Range("AY39").Value = ""

# Row 40
$ws.Range("A40").Value = 112211876
$ws.Range("B40").Value = 90689
$ws.Range("C40").Value = "Ovaliderad"
$ws.Range("D40").Value = "NT"
$ws.Range("E40").Value = 5966
$ws.Range("F40").Value = "Motaggsvamp"
$ws.Range("G40").Value = "Sarcodon squamosus"
$ws.Range("H40").Value = "(Schaeff.) Quél."
$ws.Range("I40").NumberFormat = "@"
$ws.Range("I40").Value = "1"
$ws.Range("J40").Value = "fruktkroppar"
$ws.Range("K40").NumberFormat = "@"
$ws.Range("K40").Value = ""
$ws.Range("P40").Value = "SO Simsbodarna, Dlr"
$ws.Range("Q40").Value = 515371
$ws.Range("R40").Value = 6704633
$ws.Range("S40").Value = 1
$ws.Range("T40").Value = "Dalarna"
$ws.Range("U40").Value = "Borlänge"
$ws.Range("V40").Value = "Dalarna"
$ws.Range("W40").Value = "Stora Tuna"
$ws.Range("Y40").NumberFormat = "@"
$ws.Range("Y40").Value = "2023-09-20"
$ws.Range("Z40").Value = "11:43"
$ws.Range("AA40").NumberFormat = "@"
$ws.Range("AA40").Value = "2023-09-20"
$ws.Range("AB40").Value = "11:43"
$ws.Range("AD40").Value = $false
$ws.Range("AE40").Value = $false
$ws.Range("AG40").Value = $false
$ws.Range("AT40").NumberFormat = "@"
$ws.Range("AT40").Value = ""
$ws.Range("AW40").Value = "Lars-Erik Nilsson"
$ws.Range("AX40").Value = "Lars-Erik Nilsson"
$ws.Range("AY40").NumberFormat = "@"
$ws.Range("AY40").Value = ""

# Row 41
$ws.Range("A41").Value = 112212286
$ws.Range("B41").Value = 90689
$ws.Range("C41").Value = "Ovaliderad"
$ws.Range("D41").Value = "NT"
$ws.Range("E41").Value = 5966
$ws.Range("F41").Value = "Motaggsvamp"
$ws.Range("G41").Value = "Sarcodon squamosus"
$ws.Range("H41").Value = "(Schaeff.) Quél."
$ws.Range("I41").NumberFormat = "@"
$ws.Range("I41").Value = "3"
$ws.Range("J41").Value = "fruktkroppar"
$ws.Range("K41").NumberFormat = "@"
$ws.Range("K41").Value = ""
$ws.Range("P41").Value = "Simsbodarna SO, Dlr"
$ws.Range("Q41").Value = 515476
$ws.Range("R41").Value = 6704606
$ws.Range("S41").Value = 6
$ws.Range("T41").Value = "Dalarna"
$ws.Range("U41").Value = "Borlänge"
$ws.Range("V41").Value = "Dalarna"
$ws.Range("W41").Value = "Stora Tuna"
$ws.Range("Y41").NumberFormat = "@"
$ws.Range("Y41").Value = "2023-09-20"
$ws.Range("Z41").Value = "12:08"
$ws.Range("AA41").NumberFormat = "@"
$ws.Range("AA41").Value = "2023-09-20"
$ws.Range("AB41").Value = "12:08"
$ws.Range("AD41").Value = $false
$ws.Range("AE41").Value = $false
$ws.Range("AG41").Value = $false
$ws.Range("AT41").NumberFormat = "@"
$ws.Range("AT41").Value = ""
$ws.Range("AW41").Value = "Lars-Erik Nilsson"
$ws.Range("AX41").Value = "Lars-Erik Nilsson"
$ws.Range("AY41").NumberFormat = "@"
$ws.Range("AY41").Value = ""
